$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.869.85'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '1.729.32'
$ws.Range('D4').Value = '''0.9979'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '''241.76'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').Value = '''0.9983'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '''0.4898'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '''0.2600'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').Value = '''0.06220'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('D10').Value = '1.733.43'
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').Value = '''16.03'
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('D12').Value = '''0.06915'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').Value = '''0.6099'
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').Value = '''4.496'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '''77.35'
$ws.Range('D16').Value = '''0.9986'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '26.642.43'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').Value = '''0.000007188'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').Value = '1.957.20'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('D22').Value = '''4.435'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '''8.571'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').Value = '''5.128'
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').Value = '''138.57'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').Value = '''15.34'
$ws.Range('D27').Value = '''1.783'
$ws.Range('E27').Value = '  +4.98%  '
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('D29').Value = '''106.32'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').Value = '''3.957'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').Value = '''0.07998'
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').Value = '''3.690'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').Value = '''0.04535'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''2.610'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''1.009'
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''0.6250'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '''0.9336'
$ws.Range('E37').Value = '  +2.63%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''2.058'
$ws.Range('E38').Value = '  +5.51%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''2.440'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '''0.9993'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.01504'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''5.665'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''99.76'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.3866'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '''6.946'
$ws.Range('E45').Value = '  +3.65%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '''0.1163'
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.05391'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''7.952'
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '''30.23'
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''51.74'
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').Value = '''1.238'
$ws.Range('E51').Value = '  +0.10%  '
